$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents() | Out-Null

$ws.Range("H33").Value = 266.90323
$ws.Range("I33").Value = 204.24
$ws.Range("J33").Value = 528
$ws.Range("K33").Value = 204.24
$ws.Range("L33").Value = 528
$ws.Range("M33").Value = 24.75999999999999
$ws.Range("N33").Value = -986

$ws.Range("H38").Value = 1564.439
$ws.Range("I38").Value = 188.3
$ws.Range("J38").Value = 2008.3549
$ws.Range("K38").Value = 564.9000000000001
$ws.Range("L38").Value = 6025.0647
$ws.Range("M38").Value = -192.9000000000001
$ws.Range("N38").Value = -6769.0647

$ws.Range("H58").Value = 1687.9
$ws.Range("J58").Value = 3537
$ws.Range("L58").Value = 10611
$ws.Range("N58").Value = -10911

$ws.Range("H74").Value = 3500.75
$ws.Range("I74").Value = 3467.6667
$ws.Range("J74").Value = 3600
$ws.Range("K74").Value = 3467.6667
$ws.Range("L74").Value = 3600
$ws.Range("M74").Value = -2531.6667
$ws.Range("N74").Value = -5472

$ws.Range("H77").Value = 3500.75
$ws.Range("I77").Value = 3467.6667
$ws.Range("J77").Value = 3600
$ws.Range("K77").Value = 17338.3335
$ws.Range("L77").Value = 18000
$ws.Range("M77").Value = -12658.3335
$ws.Range("N77").Value = -27360

$ws.Range("H107").Value = 2842.5715
$ws.Range("I107").Value = 3145.4546
$ws.Range("J107").Value = 2509.4
$ws.Range("K107").Value = 3145.4546
$ws.Range("L107").Value = 2509.4
$ws.Range("M107").Value = -1225.4546
$ws.Range("N107").Value = -6349.4

$ws.Range("H112").Value = 3788.2222
$ws.Range("J112").Value = 4124.25
$ws.Range("L112").Value = 12372.75
$ws.Range("N112").Value = -14588.75

$ws.Range("H116").Value = 3232.5454
$ws.Range("I116").Value = 2278.6667
$ws.Range("K116").Value = 2278.6667
$ws.Range("M116").Value = 1163.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 878.13336
$ws.Range("I2").Value = 687.75
$ws.Range("K2").Value = 687.75
$ws.Range("M2").Value = -574.75

$ws.Range("H32").Value = 5216.5405
$ws.Range("I32").Value = 5216.5405
$ws.Range("K32").Value = 5216.5405
$ws.Range("M32").Value = -4929.5405

$ws.Range("H97").Value = 805
$ws.Range("I97").Value = 540
$ws.Range("J97").Value = 1600
$ws.Range("K97").Value = 540
$ws.Range("L97").Value = 1600
$ws.Range("M97").Value = -44
$ws.Range("N97").Value = -2592

$ws.Range("H101").Value = 34666.332
$ws.Range("J101").Value = 34666.332
$ws.Range("L101").Value = 34666.332
$ws.Range("N101").Value = -41156.332

$ws.Range("H116").Value = 878.13336
$ws.Range("I116").Value = 687.75
$ws.Range("K116").Value = 687.75
$ws.Range("M116").Value = 1606.25

$ws.Range("H133").Value = 34563
$ws.Range("J133").Value = 34563
$ws.Range("L133").Value = 34563
$ws.Range("N133").Value = -39623

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 878.13336
$ws.Range("I3").Value = 687.75
$ws.Range("K3").Value = 687.75
$ws.Range("M3").Value = -573.75

$ws.Range("H132").Value = 1458045
$ws.Range("J132").Value = 1458045
$ws.Range("L132").Value = 1458045
$ws.Range("N132").Value = -1468165

$ws.Range("H134").Value = 6609.6
$ws.Range("I134").Value = 1079.6666
$ws.Range("J134").Value = 23199.4
$ws.Range("K134").Value = 3238.9998
$ws.Range("L134").Value = 69598.20000000001
$ws.Range("M134").Value = -703.9998000000001
$ws.Range("N134").Value = -74668.20000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 116815.336
$ws.Range("I22").Value = 148
$ws.Range("J22").Value = 350150
$ws.Range("K22").Value = 148
$ws.Range("L22").Value = 350150
$ws.Range("M22").Value = 202
$ws.Range("N22").Value = -350850

$ws.Range("H86").Value = 3948610
$ws.Range("J86").Value = 28813.25
$ws.Range("L86").Value = 28813.25
$ws.Range("N86").Value = -31059.25

$ws.Range("H89").Value = 3948610
$ws.Range("J89").Value = 28813.25
$ws.Range("L89").Value = 144066.25
$ws.Range("N89").Value = -155298.25

$ws.Range("H132").Value = 3182.2
$ws.Range("I132").Value = 2764.8
$ws.Range("J132").Value = 3599.6
$ws.Range("K132").Value = 8294.400000000001
$ws.Range("L132").Value = 10798.8
$ws.Range("M132").Value = -5764.400000000001
$ws.Range("N132").Value = -15858.8

$ws.Range("H141").Value = 565312.75
$ws.Range("J141").Value = 565312.75
$ws.Range("L141").Value = 565312.75
$ws.Range("N141").Value = -575672.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1239.7391
$ws.Range("I5").Value = 1416.8889
$ws.Range("K5").Value = 4250.6667
$ws.Range("M5").Value = -4138.6667

$ws.Range("H26").Value = 260
$ws.Range("J26").Value = 300
$ws.Range("L26").Value = 900
$ws.Range("N26").Value = -1476

$ws.Range("H31").Value = 1825.75
$ws.Range("I31").Value = 650.5
$ws.Range("K31").Value = 1951.5
$ws.Range("M31").Value = -1663.5

$ws.Range("H58").Value = 3200
$ws.Range("J58").Value = 3200
$ws.Range("L58").Value = 9600
$ws.Range("N58").Value = -9856

$ws.Range("H87").Value = 3019
$ws.Range("J87").Value = 3925
$ws.Range("L87").Value = 11775
$ws.Range("N87").Value = -14271

$ws.Range("H90").Value = 3019
$ws.Range("J90").Value = 3925
$ws.Range("L90").Value = 35325
$ws.Range("N90").Value = -47805

$ws.Range("H115").Value = 4936
$ws.Range("J115").Value = 5680
$ws.Range("L115").Value = 17040
$ws.Range("N115").Value = -19390

$ws.Range("H121").Value = 1218.3334
$ws.Range("J121").Value = 1402
$ws.Range("L121").Value = 4206
$ws.Range("N121").Value = -6826

$ws.Range("H131").Value = 14928257
$ws.Range("J131").Value = 3381.375
$ws.Range("L131").Value = 10144.125
$ws.Range("N131").Value = -20224.125

$ws.Range("H135").Value = 1239.7391
$ws.Range("I135").Value = 1416.8889
$ws.Range("K135").Value = 12752.0001
$ws.Range("M135").Value = -10217.0001

$ws.Range("H140").Value = 23821.146
$ws.Range("J140").Value = 3456.7144
$ws.Range("L140").Value = 10370.1432
$ws.Range("N140").Value = -20730.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2167.3333
$ws.Range("I126").Value = 1809.1111
$ws.Range("K126").Value = 5427.3333
$ws.Range("M126").Value = -2957.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents() | Out-Null

$ws.Range("H22").Value = 715.61536
$ws.Range("I22").Value = 476.16666
$ws.Range("K22").Value = 476.16666
$ws.Range("M22").Value = -181.16666

$ws.Range("H27").Value = 715.61536
$ws.Range("I27").Value = 476.16666
$ws.Range("K27").Value = 476.16666
$ws.Range("M27").Value = -369.16666

$ws.Range("H46").Value = 3648.7144
$ws.Range("I46").Value = 596.6667
$ws.Range("K46").Value = 596.6667
$ws.Range("M46").Value = -408.6667

$ws.Range("H55").Value = 284.75
$ws.Range("I55").Value = 246.91667
$ws.Range("J55").Value = 341.5
$ws.Range("K55").Value = 246.91667
$ws.Range("L55").Value = 341.5
$ws.Range("M55").Value = -73.91667000000001
$ws.Range("N55").Value = -687.5

$ws.Range("H122").Value = 35716000
$ws.Range("I122").Value = 50001616
$ws.Range("J122").Value = 1952.5
$ws.Range("K122").Value = 150004848
$ws.Range("L122").Value = 5857.5
$ws.Range("M122").Value = -150002398
$ws.Range("N122").Value = -10757.5

$ws.Range("H123").Value = 40954
$ws.Range("J123").Value = 40954
$ws.Range("L123").Value = 40954
$ws.Range("N123").Value = -50754

$ws.Range("H133").Value = 40145
$ws.Range("J133").Value = 40145
$ws.Range("L133").Value = 40145
$ws.Range("N133").Value = -45205

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents() | Out-Null

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents() | Out-Null

$ws.Range("H122").Value = 19231944
$ws.Range("I122").Value = 20834538
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 62503614
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -62501164
$ws.Range("N122").Value = -7300

$ws.Range("H136").Value = 1208.1482
$ws.Range("I136").Value = 1295.5
$ws.Range("J136").Value = 1081.091
$ws.Range("K136").Value = 3886.5
$ws.Range("L136").Value = 3243.273
$ws.Range("M136").Value = -1336.5
$ws.Range("N136").Value = -8343.272999999999
